# Auto-generated edit script: update crypto price table per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '37.268.91'
$cell.ClearFormats()

# Row 3
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.066.16'
$cell.ClearFormats()
$ws.Cells.Item(3, 5).Value = '  +1.38%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.03%  '

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '249.62'
$cell.ClearFormats()
$ws.Cells.Item(5, 5).Value = '  +0.72%  '

# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.667'
$cell.ClearFormats()
$ws.Cells.Item(6, 5).Value = '  +0.80%  '

# Row 7
$ws.Cells.Item(7, 2).Value = 'Solana'
$ws.Cells.Item(7, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '58.87'
$cell.ClearFormats()
$ws.Cells.Item(7, 5).Value = '  +3.19%  '

# Row 8
$ws.Cells.Item(8, 2).Value = 'USDC'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.ClearFormats()
$ws.Cells.Item(8, 5).Value = '  -0.07%  '

# Row 9
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.387'
$cell.ClearFormats()
$ws.Cells.Item(9, 5).Value = '  +1.87%  '

# Row 10
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0792'
$cell.ClearFormats()
$ws.Cells.Item(10, 5).Value = '  +2.20%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +1.42%  '

# Row 12
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '16.16'
$cell.ClearFormats()
$ws.Cells.Item(12, 5).Value = '  +1.77%  '

# Row 13
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.915'
$cell.ClearFormats()
$ws.Cells.Item(13, 5).Value = '  +15.62%  '

# Row 14
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.367.36'
$cell.ClearFormats()
$ws.Cells.Item(14, 5).Value = '  +1.43%  '

# Row 15
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.80'
$cell.ClearFormats()
$ws.Cells.Item(15, 5).Value = '  +4.00%  '

# Row 16
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.071.16'
$cell.ClearFormats()
$ws.Cells.Item(16, 5).Value = '  +1.96%  '

# Row 17
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '19.23'
$cell.ClearFormats()
$ws.Cells.Item(17, 5).Value = '  +16.50%  '

# Row 18
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '37.222.14'
$cell.ClearFormats()
$ws.Cells.Item(18, 5).Value = '  +0.86%  '

# Row 19
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '75.77'
$cell.ClearFormats()
$ws.Cells.Item(19, 5).Value = '  +2.36%  '

# Row 20
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0907'
$cell.ClearFormats()
$ws.Cells.Item(20, 5).Value = '  +1.56%  '

# Row 21
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.50'
$cell.ClearFormats()
$ws.Cells.Item(21, 5).Value = '  +3.71%  '

# Row 22
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '239.40'
$cell.ClearFormats()
$ws.Cells.Item(22, 5).Value = '  +1.73%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  -0.03%  '

# Row 24
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.49'
$cell.ClearFormats()
$ws.Cells.Item(24, 5).Value = '  +5.50%  '

# Row 25
$ws.Cells.Item(25, 2).Value = 'Cosmos'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.64'
$cell.ClearFormats()
$ws.Cells.Item(25, 5).Value = '  +5.90%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'PancakeSwap'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.22'
$cell.ClearFormats()
$ws.Cells.Item(26, 5).Value = '  -1.22%  '

# Row 27
$ws.Cells.Item(27, 2).Value = 'Monero'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '172.12'
$cell.ClearFormats()
$ws.Cells.Item(27, 5).Value = '  +2.23%  '

# Row 28
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '20.45'
$cell.ClearFormats()
$ws.Cells.Item(28, 5).Value = '  +3.61%  '

# Row 29
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.126'
$cell.ClearFormats()
$ws.Cells.Item(29, 5).Value = '  +1.39%  '

# Row 30
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.36'
$cell.ClearFormats()
$ws.Cells.Item(30, 5).Value = '  +14.73%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +3.21%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'Hedera'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0630'
$cell.ClearFormats()
$ws.Cells.Item(32, 5).Value = '  +3.30%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.71'
$cell.ClearFormats()
$ws.Cells.Item(33, 5).Value = '  +6.92%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.34'
$cell.ClearFormats()
$ws.Cells.Item(34, 5).Value = '  +6.54%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'Kaspa'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0881'
$cell.ClearFormats()
$ws.Cells.Item(35, 5).Value = '  +2.23%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'BinanceUSD'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.ClearFormats()
$ws.Cells.Item(36, 5).Value = '  -0.04%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'WEMIXToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.85'
$cell.ClearFormats()
$ws.Cells.Item(37, 5).Value = '  +5.61%  '

# Row 38
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.35'
$cell.ClearFormats()
$ws.Cells.Item(38, 5).Value = '  +0.85%  '

# Row 39
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.17'
$cell.ClearFormats()
$ws.Cells.Item(39, 5).Value = '  +3.46%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'Cronos'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.103'
$cell.ClearFormats()
$ws.Cells.Item(40, 5).Value = '  -2.61%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'THORChain'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.17'
$cell.ClearFormats()
$ws.Cells.Item(41, 5).Value = '  +6.97%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Aave'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '102.64'
$cell.ClearFormats()
$ws.Cells.Item(42, 5).Value = '  +7.91%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'VeChain'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0227'
$cell.ClearFormats()
$ws.Cells.Item(43, 5).Value = '  +3.20%  '

# Row 44
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.18'
$cell.ClearFormats()
$ws.Cells.Item(44, 5).Value = '  +6.23%  '

# Row 45
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '17.39'
$cell.ClearFormats()
$ws.Cells.Item(45, 5).Value = '  +0.13%  '

# Row 46
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.46'
$cell.ClearFormats()
$ws.Cells.Item(46, 5).Value = '  +1.38%  '

# Row 47
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.317.22'
$cell.ClearFormats()
$ws.Cells.Item(47, 5).Value = '  +3.50%  '

# Row 48
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.89'
$cell.ClearFormats()
$ws.Cells.Item(48, 5).Value = '  +1.76%  '

# Row 49
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.77'
$cell.ClearFormats()
$ws.Cells.Item(49, 5).Value = '  +25.63%  '

# Row 50
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.94'
$cell.ClearFormats()
$ws.Cells.Item(50, 5).Value = '  +3.96%  '

# Row 51
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.252.01'
$cell.ClearFormats()
$ws.Cells.Item(51, 5).Value = '  +1.43%  '

